$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - update column F values
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F11").Value = 7662
$wsExpo.Range("F25").Value = 290
$wsExpo.Range("F26").Value = 3650
$wsExpo.Range("F30").Value = 264
$wsExpo.Range("F35").Value = 1475
$wsExpo.Range("F46").Value = 228

# Sheet "全部类型" (all types) - same events duplicated, update column F values
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F13").Value = 7662
$wsAll.Range("F26").Value = 290
$wsAll.Range("F27").Value = 3650
$wsAll.Range("F32").Value = 264
$wsAll.Range("F36").Value = 1475
$wsAll.Range("F48").Value = 228
